$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Swap the 400mAh battery BOM line for the 110mAh battery (cheaper part + new
# order link). Only the displayed text/values change - the existing
# hyperlink relationship on C3 is left alone, matching the source edit.
$ws.Cells.Item(3, 3).Value = "https://www.sparkfun.com/products/731"
$ws.Cells.Item(3, 1).Value = "Polymer Lithium Ion Battery - 110mAh"
$ws.Cells.Item(3, 4).Value = 4.95

# Move the active selection to where the user last clicked.
$ws.Range("C22").Select()
